# Apply the "Self Learning TC02 Version 2.0" update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Self Learning TC2")

# Update existing text strings
$ws.Range("A6").Value = "You know the IP Address for the website you wish to scan"
$ws.Range("A11").Value = "AC4TC2: Ensure you're able to scan the website fully from the Nmap desktop application"

# Add the new step 4 text to B14
$ws.Range("B14").Value = "4. Scroll down the information to find the active open ports. Take note of which ports are open as these may be vulnerable to attack. There should also be detailed information on a tracerout showing you where the packet information travels from and how fast it takes to do so. "

# Adjust row height for row 14 to fit the new wrapped text
$ws.Rows.Item(14).RowHeight = 102

# Update the view: scroll position and selection
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 9
